$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): swap placeholder numbers for real column titles.
$ws.Range("A1").Value = "IP_NAS_AP"
$ws.Range("B1").Value = "Inicio_de_Conexión_Dia"
$ws.Range("C1").Value = "FIN_de_Conexión_Dia"
$ws.Range("D1").Value = "Input_Octects"
$ws.Range("E1").Value = "Output_Octects"

# --- Data rows (2-9): the full rebuilt connection log.
$rows = @(
    @("192.168.247.11", "2019-02-07", "2019-03-13", "39517", "505219"),
    @("192.168.247.12", "2019-02-26", "2019-02-26", "3084", "1344"),
    @("192.168.247.11", "2019-02-11", "2019-03-14", "1170", "495"),
    @("192.168.247.19", "2019-02-19", "2019-02-19", "2554", "334"),
    @("192.168.247.11", "2019-02-12", "2019-03-06", "544518", "1881843"),
    @("192.168.247.15", "2019-03-28", "2019-03-28", "16893", "23855"),
    @("192.168.247.16", "2019-03-15", "2019-05-02", "0", "0"),
    @("192.168.247.18", "2019-06-24", "2019-06-24", "0", "0")
)

$r = 2
foreach ($rec in $rows) {
    for ($col = 1; $col -le 5; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        # Force text BEFORE writing so date-like / numeric-like strings
        # (connection dates, octet counters) are kept verbatim instead of
        # being auto-coerced into real numbers/dates by Excel's parser.
        $cell.NumberFormat = "@"
        $cell.Value = $rec[$col - 1]
        # Put the cell's style back to the workbook default so the data
        # rows end up with no explicit style index, matching the
        # unstyled data cells from the original sheet.
        $cell.Style = "Normal"
    }
    $r++
}
